$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190364599227905
$ws.Range("B1").Value = 2.241583108901978
$ws.Range("C1").Value = 6.553914070129395
$ws.Range("D1").Value = 2.303236722946167
$ws.Range("E1").Value = 1.189778685569763
